$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Player moves left on canvas"

$ws.Range("A3").Value = "moving right"
$ws.Range("B3").Value = "right arrow key"
$ws.Range("C3").Value = "playerXpos increases and player image moves left on canvas"
$ws.Range("D3").Value = "Player moves right on canvas"

$ws.Range("A4").Value = "Non-Valid username"
$ws.Range("B4").Value = "no username"
$ws.Range("C4").Value = "asked for username again"
$ws.Range("D4").Value = "username prompt appears again"

$ws.Range("A5").Value = "Valid username"
$ws.Range("B5").Value = "username entered"
$ws.Range("C5").Value = "prompt only appears once"
$ws.Range("D5").Value = "prompt only appears once"

$ws.Range("A6").Value = "loser screen"
$ws.Range("B6").Value = "lose all lives"
$ws.Range("C6").Value = "loser screen appears"
$ws.Range("D6").Value = "loser screen appears"

$ws.Range("A7").Value = "winner screen"
$ws.Range("B7").Value = "survive 15 seconds"
$ws.Range("C7").Value = "winner screen appears"
$ws.Range("D7").Value = "winner screen appears"

$ws.Range("A8").Value = "Easy mode"
$ws.Range("B8").Value = "1"
$ws.Range("C8").Value = "Least lives and least jellyfish"
$ws.Range("D8").Value = "Least lives and least jellyfish"

$ws.Range("A9").Value = "Medium mode"
$ws.Range("B9").Value = "2"
$ws.Range("C9").Value = "Middle amount of lives and middle amount of jellyfish"
$ws.Range("D9").Value = "Middle amount of lives and middle amount of jellyfish"

$ws.Range("A10").Value = "Hard mode"
$ws.Range("B10").Value = "3"
$ws.Range("C10").Value = "Most lives and Most jellyfish"
$ws.Range("D10").Value = "Most lives and Most jellyfish"

$ws.Range("A11").Value = "Collosions"
$ws.Range("B11").Value = "Player gets struck by enemy (jellyfish)"
$ws.Range("C11").Value = "Lives go down"
$ws.Range("D11").Value = "Lives go down"

$ws.Range("A12").Value = "Lives"
$ws.Range("B12").Value = "Lives get to 0"
$ws.Range("C12").Value = "Game stops and loser screen appears"
$ws.Range("D12").Value = "Game stops and loser screen appears"

$ws.Range("A13").Value = "Time to beat"
$ws.Range("B13").Value = "survive the time to beat"
$ws.Range("C13").Value = "The winner screen appear after time to survive"
$ws.Range("D13").Value = "The winner screen appear after time to survive"

$ws.Range("A14").Value = "Player boundarys"
$ws.Range("B14").Value = "Player hits wall"
$ws.Range("C14").Value = "Player cannot progress further in that direction"
$ws.Range("D14").Value = "Player cannot progress further in that direction"

$ws.Range("A15").Value = "Wrong number (difficulty)"
$ws.Range("B15").Value = "Any other number than 1, 2 or 3"
$ws.Range("C15").Value = "Nothing happens"
$ws.Range("D15").Value = "Nothing happens"

$ws.Range("A16").Value = "Cancel or nothing in prompt"
$ws.Range("B16").Value = "cancel button or nothing in propmt"
$ws.Range("C16").Value = "Prompt reappears"
$ws.Range("D16").Value = "Prompt reappears"

$ws.Range("A17").Select()
